$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "photo" columns (front photo / ingredients photo / nutrition-facts
# photo links) were mistakenly placed in columns V, W, X. They belong
# three columns further right, in Y, Z, AA (columns that were
# previously just empty placeholder cells). This moves the header
# (row 2), the two sample rows (3 and 4), their cell comments (row 2)
# and their hyperlinks (rows 3 and 4) accordingly.

$sourceCols = @("V", "W", "X")
$destCols   = @("Y", "Z", "AA")

# ---------------------------------------------------------------------
# 1) Move the 3 cell comments anchored on row 2 (V2/W2/X2 -> Y2/Z2/AA2).
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $sourceCols.Length; $i++) {
    $srcAddr = "$($sourceCols[$i])2"
    $dstAddr = "$($destCols[$i])2"

    $comment = $ws.Range($srcAddr).Comment
    if ($comment) {
        $commentText = $comment.Text()
        $comment.Delete()
        $ws.Range($dstAddr).AddComment($commentText) | Out-Null
    }
}

# ---------------------------------------------------------------------
# 2) Move the values/styles of rows 2, 3 and 4 from V:X to Y:AA, and
#    clear out the old source cells so they no longer exist.
# ---------------------------------------------------------------------
foreach ($row in 2, 3, 4) {
    for ($i = 0; $i -lt $sourceCols.Length; $i++) {
        $srcAddr = "$($sourceCols[$i])$row"
        $dstAddr = "$($destCols[$i])$row"
        $ws.Range($srcAddr).Copy($ws.Range($dstAddr))
    }
}
foreach ($row in 2, 3, 4) {
    for ($i = 0; $i -lt $sourceCols.Length; $i++) {
        $ws.Range("$($sourceCols[$i])$row").Clear()
    }
}

# ---------------------------------------------------------------------
# 3) Re-create the hyperlinks (rows 3 and 4) on their new cells. The
#    interop's Range.Hyperlinks.Delete() removes every hyperlink on the
#    whole sheet (not just the target range), so the existing ones are
#    wiped first and then all six are re-added at their correct,
#    shifted locations with the same target URLs as before.
# ---------------------------------------------------------------------
$hyperlinkUrls = @(
    "https://drive.google.com/file/d/1N6K0GOy8ZrO21732XF1M6klhRrzyAaaP/view?usp=sharing",
    "https://drive.google.com/file/d/14181RIHe89KN2aPzT3C2NQYqFZficOwE/view?usp=drive_link",
    "https://drive.google.com/file/d/1e8sW-Xp5QXc41V5ry-hLn8T8s13fvS38/view?usp=sharing"
)

if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Range("A1").Hyperlinks.Delete()
}

$tmp = $ws.Range("ZZ1")
foreach ($row in 3, 4) {
    for ($i = 0; $i -lt $destCols.Length; $i++) {
        $dst = $ws.Range("$($destCols[$i])$row")

        # Stash the current formatting so it can be restored after
        # Hyperlinks.Add() forces the built-in "Hyperlink" style on it.
        $dst.Copy($tmp)

        $ws.Hyperlinks.Add($dst, $hyperlinkUrls[$i]) | Out-Null

        $tmp.Copy()
        $dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    }
}
$tmp.Clear()
$excel.CutCopyMode = $false
